$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("CreateUser")
$ws.Range("B2").Value = "Wed Apr 30 22:41:30 IST 2025"
$ws.Range("B3").Value = "Wed Apr 30 22:42:16 IST 2025"
$ws.Range("B4").Value = "Wed Apr 30 22:42:50 IST 2025"
$ws = $wb.Worksheets.Item("CreateUserSpChar")
$ws.Range("B2").Value = "Wed Apr 30 22:43:19 IST 2025"
$ws.Range("B3").Value = "Wed Apr 30 22:43:46 IST 2025"
$ws.Range("B4").Value = "Wed Apr 30 22:44:14 IST 2025"
$ws = $wb.Worksheets.Item("CreateUserSpCharError")
$ws.Range("B2").Value = "Wed Apr 30 22:44:42 IST 2025"
$ws.Range("B3").Value = "Wed Apr 30 22:45:01 IST 2025"
$ws.Range("B4").Value = "Wed Apr 30 22:45:20 IST 2025"
$ws.Range("B5").Value = "Wed Apr 30 22:45:40 IST 2025"
$ws.Range("B6").Value = "Wed Apr 30 22:45:59 IST 2025"
$ws.Range("B7").Value = "Wed Apr 30 22:46:18 IST 2025"
$ws = $wb.Worksheets.Item("FindUser")
$ws.Range("B2").Value = "Wed Apr 30 22:46:37 IST 2025"
$ws.Range("B3").Value = "Wed Apr 30 22:46:59 IST 2025"
$ws.Range("B4").Value = "Wed Apr 30 22:47:20 IST 2025"
$ws.Range("B5").Value = "Wed Apr 30 22:47:42 IST 2025"
$ws.Range("B6").Value = "Wed Apr 30 22:48:05 IST 2025"
$ws.Range("B7").Value = "Wed Apr 30 22:48:28 IST 2025"
$ws.Range("B8").Value = "Wed Apr 30 22:48:52 IST 2025"
$ws.Range("B9").Value = "Wed Apr 30 22:49:14 IST 2025"
$ws.Range("B10").Value = "Wed Apr 30 22:49:35 IST 2025"
$ws.Range("B11").Value = "Wed Apr 30 22:49:59 IST 2025"
$ws.Range("B12").Value = "Wed Apr 30 22:50:20 IST 2025"
$ws.Range("B13").Value = "Wed Apr 30 22:50:43 IST 2025"
$ws = $wb.Worksheets.Item("PassCase")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Wed Apr 30 22:51:08 IST 2025"
$ws = $wb.Worksheets.Item("UsernameCase")
$ws.Range("B2").Value = "Wed Apr 30 22:51:20 IST 2025"
$ws.Range("B3").Value = "Wed Apr 30 22:51:36 IST 2025"
$ws.Range("B4").Value = "Wed Apr 30 22:51:52 IST 2025"
$ws = $wb.Worksheets.Item("CreateUserPasswordSpChar")
$ws.Range("B2").Value = "Wed Apr 30 22:52:10 IST 2025"
$ws.Range("B3").Value = "Wed Apr 30 22:52:37 IST 2025"
$ws.Range("B4").Value = "Wed Apr 30 22:53:05 IST 2025"
$ws.Range("B5").Value = "Wed Apr 30 22:53:33 IST 2025"
$ws.Range("B6").Value = "Wed Apr 30 22:53:59 IST 2025"
$ws.Range("B7").Value = "Wed Apr 30 22:54:29 IST 2025"
$ws.Range("B8").Value = "Wed Apr 30 22:54:55 IST 2025"
$ws.Range("B9").Value = "Wed Apr 30 22:55:23 IST 2025"
$ws.Range("B10").Value = "Wed Apr 30 22:55:51 IST 2025"
$ws.Range("B11").Value = "Wed Apr 30 22:56:22 IST 2025"
$ws.Range("B12").Value = "Wed Apr 30 22:56:51 IST 2025"
$ws.Range("B13").Value = "Wed Apr 30 22:57:18 IST 2025"
$ws.Range("B14").Value = "Wed Apr 30 22:57:45 IST 2025"
$ws.Range("B15").Value = "Wed Apr 30 22:58:12 IST 2025"
$ws.Range("B16").Value = "Wed Apr 30 22:58:40 IST 2025"
$ws.Range("B17").Value = "Wed Apr 30 22:59:09 IST 2025"
$ws = $wb.Worksheets.Item("ModifyUserPwd")
$ws.Range("B2").Value = "Wed Apr 30 23:01:02 IST 2025"
$ws.Range("B3").Value = "Wed Apr 30 23:01:27 IST 2025"
$ws.Range("B4").Value = "Wed Apr 30 23:01:55 IST 2025"
$ws.Range("B5").Value = "Wed Apr 30 23:02:24 IST 2025"
$ws.Range("B6").Value = "Wed Apr 30 23:02:52 IST 2025"
$ws.Range("B7").Value = "Wed Apr 30 23:03:24 IST 2025"
$ws.Range("B8").Value = "Wed Apr 30 23:03:48 IST 2025"
$ws = $wb.Worksheets.Item("ModifyUser")
$ws.Range("B2").Value = "Wed Apr 30 23:04:16 IST 2025"
$ws.Range("B3").Value = "Wed Apr 30 23:04:49 IST 2025"
$ws = $wb.Worksheets.Item("CreateUserSCFNameErr")
$ws.Range("B2").Value = "Wed Apr 30 23:05:25 IST 2025"
$ws.Range("B3").Value = "Wed Apr 30 23:05:50 IST 2025"
$ws.Range("B4").Value = "Wed Apr 30 23:06:10 IST 2025"
$ws.Range("B5").Value = "Wed Apr 30 23:06:30 IST 2025"
$ws = $wb.Worksheets.Item("CreateUserSCLNameErr")
$ws.Range("B2").Value = "Wed Apr 30 23:06:53 IST 2025"
$ws.Range("B3").Value = "Wed Apr 30 23:07:12 IST 2025"
$ws.Range("B4").Value = "Wed Apr 30 23:07:31 IST 2025"
$ws.Range("B5").Value = "Wed Apr 30 23:07:55 IST 2025"
$ws = $wb.Worksheets.Item("CreateUserErrors")
$ws.Range("A13").Value = "Fail"
$ws.Range("B13").Value = "Wed Apr 30 23:08:14 IST 2025"
$ws.Range("A14").Value = "Fail"
$ws.Range("B14").Value = "Wed Apr 30 23:08:34 IST 2025"
$ws = $wb.Worksheets.Item("AddDeleteRole")
$ws.Range("B2").Value = "Wed Apr 30 23:09:53 IST 2025"
$ws.Range("B3").Value = "Wed Apr 30 23:10:24 IST 2025"
$ws.Range("B4").Value = "Wed Apr 30 23:10:54 IST 2025"
$ws.Range("B5").Value = "Wed Apr 30 23:11:19 IST 2025"
$ws = $wb.Worksheets.Item("SearchRole")
$ws.Range("B2").Value = "Wed Apr 30 23:11:45 IST 2025"
$ws.Range("B3").Value = "Wed Apr 30 23:12:13 IST 2025"
$ws.Range("B4").Value = "Wed Apr 30 23:12:34 IST 2025"
$ws.Range("B5").Value = "Wed Apr 30 23:12:56 IST 2025"
